# Fill intended time for the week (row 19/20, "Планируемые часы работы")
# and update today's actual time (row 22/23, "Фактические часы работы").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19/20: intended (planned) hours for Thu/Fri/Sat (G:L) ---
$ws.Range("G19").Value2 = 0.41666666666666669
$ws.Range("H19").Value2 = 7
$ws.Range("I19").Value2 = 0.625
$ws.Range("J19").Value2 = 3.5
$ws.Range("K19").Value2 = 0.5
$ws.Range("L19").Value2 = 5

$ws.Range("G20").Value2 = 0.70833333333333337
$ws.Range("I20").Value2 = 0.77083333333333337
$ws.Range("K20").Value2 = 0.70833333333333337

# --- Row 22/23: today's actual hours (E:F) ---
$ws.Range("E22").Value2 = 0.59027777777777779
$ws.Range("F22").Value2 = 5.5

$ws.Range("E23").Value2 = 0.81944444444444453

# --- Update the view: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2

$ws.Range("K22").Select()
